$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overlay")

$ws.Range("B7").Value = "Zero"
$ws.Range("A7").Value = "6"
$ws.Range("C7").Value = "<Type=Text><Text=0.00>"
$ws.Range("D7").Value = "<X=360><Y=442>"
$ws.Range("D7").NumberFormat = "@"
